# "Rename key to cd"
# The two placeholder labels that used to end in "_key" are renamed to end
# in "_cd": [project_key] -> [project_cd] and [bill_key] -> [bill_cd].
# Everything else on the sheet (ordering, other labels) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "[project_cd]"
$ws.Range("A5").Value = "[bill_cd]"
